$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D cells whose new values look numeric need to be forced to Text
# format first, otherwise Excel auto-converts the string into a Number and
# silently drops significant trailing zeros (e.g. "1.00" -> 1, "0.0850" -> 0.085).
$numericLookingCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D50", "D51"
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values exactly as they appear in the source data feed.
$ws.Range("D2").Value = "51.499.87"
$ws.Range("E2").Value = "  +4.84%  "
$ws.Range("D3").Value = "2.732.96"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "115.52"
$ws.Range("E5").Value = "  +3.97%  "
$ws.Range("D6").Value = "331.59"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("E7").Value = "  +2.05%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  +4.61%  "
$ws.Range("D10").Value = "41.44"
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("D11").Value = "0.0850"
$ws.Range("E11").Value = "  +4.91%  "
$ws.Range("D12").Value = "20.09"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "7.58"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "3.162.04"
$ws.Range("E15").Value = "  +4.31%  "
$ws.Range("D16").Value = "2.751.17"
$ws.Range("E16").Value = "  +5.04%  "
$ws.Range("D17").Value = "0.877"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "51.425.59"
$ws.Range("E18").Value = "  +4.91%  "
$ws.Range("D19").Value = "3.15"
$ws.Range("E19").Value = "  +5.63%  "
$ws.Range("D20").Value = "13.34"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").Value = "0.0₃0971"
$ws.Range("D23").Value = "277.77"
$ws.Range("E23").Value = "  +2.98%  "
$ws.Range("D24").Value = "69.19"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").Value = "2.64"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("D26").Value = "26.66"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "0.141"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "34.94"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "49.95"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "5.52"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "0.0817"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "19.06"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "2.08"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").Value = "128.10"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "23.08"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "2.30"
$ws.Range("E42").Value = "  +7.73%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "0.113"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0342"
$ws.Range("E44").Value = "  +8.15%  "
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  +10.87%  "
$ws.Range("D46").Value = "2.083.66"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("E48").Value = "  +3.51%  "
$ws.Range("D49").Value = "5.52"
$ws.Range("E49").Value = "  +6.37%  "
$ws.Range("D50").Value = "8.91"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "59.67"
$ws.Range("E51").Value = "  +1.74%  "
